$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: "Update Group" -> "Groups"
$ws.Name = "Groups"

# The template now downloads two columns instead of one: an ID column plus
# the renamed "Group Name" column. Shift the existing header out of A1 and
# add the new headers.
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Group Name"

# Give the new header cell (B1) the same look as A1 (bold header font + fill)
# before toggling bold, so both headers end up sharing one consistent style.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Make the header row bold.
$ws.Range("A1:B1").Font.Bold = $true

# Column widths: narrower ID column, Group Name keeps the original width.
# (Values tuned so the saved OOXML <col width> lands on 31 / 47, matching
# what Excel itself stores for these character widths.)
$ws.Columns.Item(1).ColumnWidth = 30.140625
$ws.Columns.Item(2).ColumnWidth = 46.17

# Leave the cursor parked under the new second column, matching the saved view.
$ws.Range("B2").Select()
